$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.755.56'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').Value = '2.282.49'
$ws.Range('E3').Value = '  -0.15%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '124.68'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +7.09%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '266.80'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.84%  '
$ws.Range('E7').Value = '  +2.35%  '
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('E9').Value = '  +1.37%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '48.35'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.32%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0949'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.57%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '9.39'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +4.56%  '
$ws.Range('E13').Value = '  -1.02%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '15.50'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.90%  '
$ws.Range('E15').Value = '  +3.96%  '
$ws.Range('D16').Value = '2.625.23'
$ws.Range('D17').Value = '2.275.56'
$ws.Range('E17').Value = '  -1.00%  '
$ws.Range('D18').Value = '43.699.26'
$ws.Range('E18').Value = '  +0.20%  '
$ws.Range('E20').Value = '  +0.33%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '72.40'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.19%  '
$ws.Range('E22').Value = '  +0.90%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '235.55'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.92%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '9.52'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -6.71%  '
$ws.Range('E25').Value = '  -2.41%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '12.00'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +2.86%  '
$ws.Range('E27').Value = '  +1.63%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '42.28'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.72%  '
$ws.Range('E29').Value = '  -0.66%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '172.98'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.32%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '21.71'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.63%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0924'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.17%  '
$ws.Range('E34').Value = '  +0.70%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.35'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +13.45%  '
$ws.Range('E36').Value = '  +2.07%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0377'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +4.92%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.64'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.36%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.106'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.07%  '
$ws.Range('E40').Value = '  +5.24%  '
$ws.Range('E41').Value = '  -3.57%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '74.05'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.81%  '
$ws.Range('E43').Value = '  -0.95%  '
$ws.Range('E44').Value = '  -0.07%  '
$ws.Range('E45').Value = '  -1.21%  '
$ws.Range('E46').Value = '  -11.29%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '74.17'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +38.06%  '
$ws.Range('B48').Value = 'TrustWalletToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.27'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.72%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '8.59'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.90%  '
$ws.Range('E50').Value = '  +0.39%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '102.17'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.45%  '
